# Update the "DB Schema" sheet to reflect the revised database schema
# (adds user_aadhar_no / total_slots_2_wheeler / total_slots_4_wheeler /
# parking_area_in_sqft fields, widens several varchar() types, drops the
# "username" field row, and re-points several FK/type cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DB Schema")

# ---- "1. User_info" table (left, columns C:E) ----------------------------
# Row 7: user_first_name
$ws.Range("D7").Value = "varchar(100)"

# Row 8: user_last_name
$ws.Range("D8").Value = "varchar(100)"

# Row 9: user_email_id
$ws.Range("D9").Value = "varchar(100)"

# Row 10: address
$ws.Range("D10").Value = "varchar(200)"

# Row 11: used to be "username" -> now "password"
$ws.Range("C11").Value = "password"
$ws.Range("D11").Value = "varchar(100)"

# Row 12: used to be "password" -> now "user_mobile_no"
$ws.Range("C12").Value = "user_mobile_no"
$ws.Range("D12").Value = "varchar(10)"

# Row 13: used to be "user_mobile_no" -> now "user_aadhar_no"
$ws.Range("C13").Value = "user_aadhar_no"
$ws.Range("D13").Value = "varchar(20)"

# Row 14: used to be "user_adhar_no" -> now "User_Role"
$ws.Range("C14").Value = "User_Role"
$ws.Range("D14").Value = "varchar(100)"

# Row 15: used to be "User_Role" -> row now belongs solely to the right table
$ws.Range("C15:E15").Clear()

# ---- "2.Parking_space_info" table (right, columns H:J) --------------------
# Row 7: parking_name
$ws.Range("I7").Value = "varchar(100)"

# Row 8: used to be "no_of_2_wheeler" -> now "total_slots_2_wheeler"
$ws.Range("H8").Value = "total_slots_2_wheeler"

# Row 9: used to be "no_of_4_wheeler" -> now "total_slots_4_wheeler"
$ws.Range("H9").Value = "total_slots_4_wheeler"
$ws.Range("I9").Value = "int"

# Row 10: per_hr_price
$ws.Range("I10").Value = "int"

# Row 11: user_id (FK)
$ws.Range("I11").Value = "int"

# Row 12: parking_address
$ws.Range("I12").Value = "varchar(200) "

# Row 13: used to be "parking_area" -> now "parking_area_in_sqft"
$ws.Range("H13").Value = "parking_area_in_sqft"
$ws.Range("I13").Value = "varchar(10)"

# Row 14: slot_available_2_wheeler
$ws.Range("I14").Value = "int"

# Row 15: slot_available_4_wheeler
$ws.Range("I15").Value = "int"

# Row 16: parking_pincode
$ws.Range("I16").Value = "int"

# Row 17: used to be "user_adhar_no" label slot -> now "parking_addr_latitude"
$ws.Range("H17").Value = "parking_addr_latitude"
$ws.Range("I17").Value = "varchar(20)"

# Row 18: used to be "parking_addr_latitude" -> now "parking_addr_longitude"
$ws.Range("H18").Value = "parking_addr_longitude"
$ws.Range("I18").Value = "varchar(20)"

# ---- "3.Bookings_info" table (columns C:E, rows 20-30) ---------------------
# Row 22: owner_id now carries the FK key marker
$ws.Range("E22").Value = "FK"

# Row 23: customer_id type int, FK key marker
$ws.Range("D23").Value = "int"
$ws.Range("E23").Value = "FK"

# Row 24: vehicle_type
$ws.Range("D24").Value = "varchar(10) "

# Row 29: booking_status
$ws.Range("D29").Value = "varchar(20)"

# Row 30: parking_id
$ws.Range("D30").Value = "int"

# ---- Selection state, as captured in the saved file ------------------------
$ws.Range("I16").Select()
